# Fix source document that used Heading1 for normal text:
# The paragraph containing "None" (under the "Prerequisite tasks" Heading2)
# was incorrectly styled as Heading1 (with direct sz=22/szCs=22 run/paragraph
# formatting) and carried an auto-generated heading bookmark. Convert it back
# to a plain Normal paragraph and drop the now-meaningless bookmark, the way
# a later "Detailed Instructions" style paragraph already looks.

$d = $word.ActiveDocument

# Locate the bookmark Word created for the old (incorrect) Heading1 paragraph
# that just says "None", then delete the bookmark first (while we still have
# a stable name to find it by) and re-style the paragraph as Normal, which
# also clears the stray direct sz/szCs formatting that had been pinned to
# match the heading's apparent size. (Hidden bookmarks like this one, whose
# name starts with "_", are still reachable by name without having to flip
# Bookmarks.ShowHidden.)
$noneBookmark = $d.Bookmarks.Item("_aqdfz55armzg")
$noneRange = $noneBookmark.Range
$noneParagraph = $noneRange.Paragraphs.Item(1)

$noneBookmark.Delete()
$noneParagraph.Style = "Normal"
